# Update factsheets with text edits from COMM
#
# The "No. of 990 Filers w/ Gov Grants" counts (column B on the detail
# sheets, column A on "Overall") were stored as real numbers. They need to
# become literal text values that use a thousands separator (e.g. 8743 ->
# "8,743"), matching the formatting already used for the dollar/percentage
# columns on these sheets. The "County" sheet is also missing its "Total"
# summary row (row 64), which every other detail sheet already has.

$wb = $excel.ActiveWorkbook

function Convert-CellToCommaText($cell) {
    $v = $cell.Value2
    $formatted = $v.ToString("N0")
    $cell.NumberFormat = "@"
    $cell.Value = $formatted
}

function Convert-ColumnBToCommaText($ws, $firstRow, $lastRow) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        Convert-CellToCommaText $ws.Cells.Item($r, 2)
    }
}

# ---- "Overall" sheet: A2 (8743 -> "8,743") ----
$wsOverall = $wb.Worksheets.Item("Overall")
Convert-CellToCommaText $wsOverall.Cells.Item(2, 1)

# ---- "County" sheet: B2:B63 ----
$wsCounty = $wb.Worksheets.Item("County")
Convert-ColumnBToCommaText $wsCounty 2 63

# ---- "Congressional District" sheet: B2:B28 (includes Total row) ----
$wsCD = $wb.Worksheets.Item("Congressional District")
Convert-ColumnBToCommaText $wsCD 2 28

# ---- "Size" sheet: B2:B8 (includes Total row) ----
$wsSize = $wb.Worksheets.Item("Size")
Convert-ColumnBToCommaText $wsSize 2 8

# ---- "Subsector" sheet: B2:B14 (includes Total row) ----
$wsSub = $wb.Worksheets.Item("Subsector")
Convert-ColumnBToCommaText $wsSub 2 14

# ---- "County" sheet: add the missing "Total" row (row 64) ----
$totalRange = $wsCounty.Range("A64:F64")
$totalRange.NumberFormat = "@"

$wsCounty.Cells.Item(64, 1).Value = "Total"
$wsCounty.Cells.Item(64, 2).Value = "8,743"
$wsCounty.Cells.Item(64, 3).Value = "`$29,720,633,621"
$wsCounty.Cells.Item(64, 4).Value = "8.36%"
$wsCounty.Cells.Item(64, 5).Value = "-12.97%"
$wsCounty.Cells.Item(64, 6).Value = "69.53%"
